$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C (Förändrad / Changed date) for all existing data rows
#    (rows 2 through 532) from 45182 to 45184.
$ws.Range("C2:C532").Value = 45184

# Row 532 gains an explicit custom row height in the target workbook.
$ws.Rows.Item(532).RowHeight = 15

# 2) Append new row 533: A 42928-2023
$ws.Cells.Item(533, 1).Value = "A 42928-2023"           # A - Beteckning
$ws.Cells.Item(533, 2).Value = 45182                     # B - Datum
$ws.Cells.Item(533, 3).Value = 45184                     # C - Förändrad
$ws.Cells.Item(533, 4).Value = "ÖSTERGÖTLANDS LÄN"       # D - Län
$ws.Cells.Item(533, 5).Value = "MOTALA"                  # E - Kommun
$ws.Cells.Item(533, 6).Value = "Kyrkan"                  # F - Markägare
$ws.Cells.Item(533, 7).Value = 3.3                       # G - Area (ha)
$ws.Cells.Item(533, 8).Value = 0                         # H - Fridlysta
$ws.Cells.Item(533, 9).Value = 0                         # I - Signalarter
$ws.Cells.Item(533, 10).Value = 0                        # J - NT
$ws.Cells.Item(533, 11).Value = 0                        # K - VU
$ws.Cells.Item(533, 12).Value = 0                        # L - EN
$ws.Cells.Item(533, 13).Value = 0                        # M - CR
$ws.Cells.Item(533, 14).Value = 0                        # N - RE
$ws.Cells.Item(533, 15).Value = 0                        # O - Rödlistade
$ws.Cells.Item(533, 16).Value = 0                        # P - Hotade
$ws.Cells.Item(533, 17).Value = 0                        # Q - Alla arter
$ws.Cells.Item(533, 18).Value = ""                       # R - Artnamn
$ws.Rows.Item(533).RowHeight = 15

# 3) Append new row 534: A 43257-2023 (no Markägare / F column value)
$ws.Cells.Item(534, 1).Value = "A 43257-2023"            # A - Beteckning
$ws.Cells.Item(534, 2).Value = 45183                     # B - Datum
$ws.Cells.Item(534, 3).Value = 45184                     # C - Förändrad
$ws.Cells.Item(534, 4).Value = "ÖSTERGÖTLANDS LÄN"       # D - Län
$ws.Cells.Item(534, 5).Value = "MOTALA"                  # E - Kommun
$ws.Cells.Item(534, 7).Value = 1.3                       # G - Area (ha)
$ws.Cells.Item(534, 8).Value = 0                         # H - Fridlysta
$ws.Cells.Item(534, 9).Value = 0                         # I - Signalarter
$ws.Cells.Item(534, 10).Value = 0                        # J - NT
$ws.Cells.Item(534, 11).Value = 0                        # K - VU
$ws.Cells.Item(534, 12).Value = 0                        # L - EN
$ws.Cells.Item(534, 13).Value = 0                        # M - CR
$ws.Cells.Item(534, 14).Value = 0                        # N - RE
$ws.Cells.Item(534, 15).Value = 0                        # O - Rödlistade
$ws.Cells.Item(534, 16).Value = 0                        # P - Hotade
$ws.Cells.Item(534, 17).Value = 0                        # Q - Alla arter
$ws.Cells.Item(534, 18).Value = ""                       # R - Artnamn

# 4) Apply the same date format (yyyy-mm-dd) and wrap-text formatting used
#    by the rest of the sheet to the new rows' relevant cells.
$ws.Range("B533:C534").NumberFormat = "YYYY-MM-DD"
$ws.Range("R533:R534").WrapText = $true
